$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- TU10 (row 11): test now passes - "Actual Results" matches expected, Pass/Fail -> Pass ---
$ws.Range("G11").Value = "Bohater ma możliwość wskoczenia na równolegle położone podłoże do aktualnego na którym się znajduje. Może dokonać tego przenikając jedynie przez dolną sciankę podłoża. Boczna scianka podłoża blokuje skok.`n- Zmiana wymagań - boczna ścianka powinna blokować skok."
$ws.Range("H11").Value = "Pass"
$ws.Range("H11").Interior.Color = 5296274

# Row 11 grew taller to fit the longer "Actual Results" text
$ws.Rows.Item(11).RowHeight = 105

# --- TU11 (row 12): test now passes - "Actual Results" matches expected, Pass/Fail -> Pass ---
$ws.Range("G12").Value = "Po uruchomieniu gry, `nładują się tekstury zgodne ze specyfikacją"
$ws.Range("H12").Value = "Pass"
$ws.Range("H12").Interior.Color = 5296274

# Reflect the final on-screen selection from the authored edit
$ws.Range("H12").Select()
